$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 55.696362090774983
$ws.Range("C2").Value = 0.47699418615428602
$ws.Range("D2").Value = 37.94723128054742
$ws.Range("E2").Value = 0.43367471640083499
$ws.Range("F2").Value = 0.65853983660886828
$ws.Range("G2").Value = 0.6030176210603817
$ws.Range("H2").Value = 0.56632528359916501
$ws.Range("I2").Value = 0.75275452343933613

$ws.Range("B3").Value = 55.800885098409537
$ws.Range("C3").Value = 0.47788934097390934
$ws.Range("D3").Value = 38.844992393528948
$ws.Range("E3").Value = 0.4353039617405825
$ws.Range("F3").Value = 0.65977569047410534
$ws.Range("G3").Value = 0.6172838995835308
$ws.Range("H3").Value = 0.5646960382594175
$ws.Range("I3").Value = 0.75214038542848971

$ws.Range("B4").Value = 56.293313910375453
$ws.Range("C4").Value = 0.48210659451768267
$ws.Range("D4").Value = 38.528877680908238
$ws.Range("E4").Value = 0.44302075806897917
$ws.Range("F4").Value = 0.66559804542154355
$ws.Range("G4").Value = 0.6122605359400165
$ws.Range("H4").Value = 0.55697924193102089
$ws.Range("I4").Value = 0.74809783802052565
